# This script reorders the data for three groups of rows in the "Artfynd" sheet.
# Group A: rows 2,3,4 are cyclically rotated (new2=old3, new3=old4, new4=old2)
# Group B: rows 7,8 are swapped
# Group C: rows 12,14 are swapped
#
# A scratch row (101), well below the used range (A1:AY27), is used as
# temporary storage while shuffling each group. Because Excel's PasteSpecial
# only writes cells that exist in the copied source range, any destination
# cell that must become blank (because the incoming row has no content for
# that column) is cleared explicitly afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- Group A: rows 2, 3, 4 ----------
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(101).PasteSpecial()   # stash old row 2

$ws.Rows.Item(3).Copy()
$ws.Rows.Item(2).PasteSpecial()     # row 2 = old row 3

$ws.Rows.Item(4).Copy()
$ws.Rows.Item(3).PasteSpecial()     # row 3 = old row 4

$ws.Rows.Item(101).Copy()
$ws.Rows.Item(4).PasteSpecial()     # row 4 = old row 2
$ws.Rows.Item(101).Delete()

# old row 2 had no "Publik kommentar" (AC), but row 4 previously did;
# clear the stale value that PasteSpecial left behind.
$ws.Range("AC4").ClearContents()

# ---------- Group B: rows 7, 8 ----------
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(101).PasteSpecial()   # stash old row 7

$ws.Rows.Item(8).Copy()
$ws.Rows.Item(7).PasteSpecial()     # row 7 = old row 8

$ws.Rows.Item(101).Copy()
$ws.Rows.Item(8).PasteSpecial()     # row 8 = old row 7
$ws.Rows.Item(101).Delete()

# old row 8 had no "Aktivitet" (M), old row 7 had no "Publik kommentar" (AC);
# clear the stale values left behind by PasteSpecial.
$ws.Range("M7").ClearContents()
$ws.Range("AC8").ClearContents()

# ---------- Group C: rows 12, 14 ----------
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(101).PasteSpecial()   # stash old row 12

$ws.Rows.Item(14).Copy()
$ws.Rows.Item(12).PasteSpecial()    # row 12 = old row 14

$ws.Rows.Item(101).Copy()
$ws.Rows.Item(14).PasteSpecial()    # row 14 = old row 12
$ws.Rows.Item(101).Delete()

# old row 14 had no "Aktivitet" (M); clear the stale value.
$ws.Range("M12").ClearContents()
